$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Laborator 21.03.2024 - Greedy, probleme rezolvate eficient
# Mark attendance for săpt. 4 (column F) for students who attended,
# and add the new students to the roster.

# --- Mark attendance (săpt. 4 / column F) for already-listed students ---
$ws.Range("F6").Value  = $true   # Andrei Negrut
$ws.Range("F9").Value  = $true   # Bianca Nicorici
$ws.Range("F20").Value = $true   # Georgiana Galea
$ws.Range("F21").Value = $true   # Luca Bulea
$ws.Range("F25").Value = $true   # Noelia Sfrangeu
$ws.Range("F33").Value = $true   # Sebastian Pop
$ws.Range("F37").Value = $true   # Vlad Varkonyi

# --- Add the new students (present in săpt. 4 / column F) ---
$ws.Range("B38").Value = "Andrei Tig"
$ws.Range("F38").Value = $true

$ws.Range("B39").Value = "Victor Pitirici"
$ws.Range("F39").Value = $true

$ws.Range("B40").Value = "Alexandra Iovan"
$ws.Range("F40").Value = $true

$ws.Range("B41").Value = "David Nadis"
$ws.Range("F41").Value = $true

# --- Re-sort the whole roster alphabetically by name (column B) ---
$sortRange = $ws.Range("B3:S41")
$keyRange = $ws.Range("B3:B41")
$sortRange.Sort($keyRange, 1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 1)

# --- Update selection to match the author's last cursor position ---
$ws.Range("K10").Select()
